$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 26
$ws.Range("B6").Value = "Update index.py"
$ws.Range("C6").Value = "riya-morankar"
$ws.Range("D6").Value = "N/A"
$ws.Range("E6").Value = "edit1 to main"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "2025-06-17"
$ws.Range("F6").Style = "Normal"
